$d = $word.ActiveDocument

# =====================================================================
# 1) Paragraph "Set up Github - July 17-18th " gains a new, bold
#    "JULY 17th [COMPLETED]" tail (with "th" superscripted).
# =====================================================================
$p8 = $d.Paragraphs(8)
$r8 = $p8.Range

$ip1 = $d.Range($r8.End - 1, $r8.End - 1)
$ip1.InsertAfter("JULY 17")
$ip1.Font.Name = "Lucida Console"
$ip1.Font.Bold = $true

$p8b = $d.Paragraphs(8)
$r8b = $p8b.Range
$ip2 = $d.Range($r8b.End - 1, $r8b.End - 1)
$ip2.InsertAfter("th")
$ip2.Font.Name = "Lucida Console"
$ip2.Font.Bold = $true
$ip2.Font.Superscript = $true

$p8c = $d.Paragraphs(8)
$r8c = $p8c.Range
$ip3 = $d.Range($r8c.End - 1, $r8c.End - 1)
$ip3.InsertAfter(" [COMPLETED]")
$ip3.Font.Name = "Lucida Console"
$ip3.Font.Bold = $true

# =====================================================================
# 2) A brand-new sub-bullet paragraph is added right after it:
#    "Fill in Parameters" - July 19-21st
# =====================================================================
$p8d = $d.Paragraphs(8)
$r8d = $p8d.Range
$r8d.InsertParagraphAfter()

$p9 = $d.Paragraphs(9)
$r9 = $p9.Range
$ip4 = $d.Range($r9.Start, $r9.Start)
$ip4.InsertAfter([char]0x201C + "Fill in Parameters" + [char]0x201D + " " + [char]0x2013 + " July 19-21")
$ip4.Font.Name = "Lucida Console"
$ip4.Font.Bold = $false

$p9b = $d.Paragraphs(9)
$r9b = $p9b.Range
$ip5 = $d.Range($r9b.End - 1, $r9b.End - 1)
$ip5.InsertAfter("st")
$ip5.Font.Name = "Lucida Console"
$ip5.Font.Bold = $false
$ip5.Font.Superscript = $true

$p9c = $d.Paragraphs(9)
$r9c = $p9c.Range
$ip6 = $d.Range($r9c.End - 1, $r9c.End - 1)
$ip6.InsertAfter(" ")
$ip6.Font.Name = "Lucida Console"
$ip6.Font.Bold = $false

# =====================================================================
# 3) The "July 15th ... Submit Plan for Application [COMPLETED]" line
#    drops its "[" / "COMPLETED" / "]" run trio and instead gets a
#    single bold "[COMPLETED]" run placed just ahead of the _GoBack
#    bookmark.
# =====================================================================
$bm = $d.Bookmarks("_GoBack")
$bracketRange = $d.Range($bm.Start - 1, $bm.Start)
$bracketRange.Delete()

$bm2 = $d.Bookmarks("_GoBack")
$ip7 = $d.Range($bm2.Start, $bm2.Start)
$ip7.InsertBefore("[COMPLETED]")
$newRunRange = $d.Range($bm2.Start - 11, $bm2.Start)
$newRunRange.Font.Name = "Lucida Console"
$newRunRange.Font.Bold = $true

$bm3 = $d.Bookmarks("_GoBack")
$lastPara = $bm3.Range.Paragraphs(1)
$oldRange = $d.Range($bm3.End, $lastPara.Range.End - 1)
$oldRange.Delete()
